# Update the underlying attrition-rate source data on the "SS_att" sheet.
# (Workbook was re-built from the latest replication data; these are the
# refreshed raw inputs - average days-to-separation / SDs / attrition
# fraction per admin group - that the rest of the workbook (Attrition,
# SS, SS_cond_survey, ...) recalculates from.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SS_att")

$ws.Range("B2").Value = 35.404255319148938
$ws.Range("C2").Value = 32.322580645161288
$ws.Range("D2").Value = 37.324324324324323
$ws.Range("E2").Value = 35.136612021857921
$ws.Range("L2").Value = 0.46218569349467797

$ws.Range("B3").Value = 3.4834961191702076
$ws.Range("C3").Value = 2.6504216953471356
$ws.Range("D3").Value = 3.1208005093831628
$ws.Range("E3").Value = 1.7862037315142396

# Downstream sheets (Attrition, etc.) hold formulas referencing SS_att, so
# they are picked up automatically by the recalculation that follows this
# script. Mark the workbook for a full recalculation on next open as well,
# matching the refreshed/rebuilt data (calcPr fullCalcOnLoad).
$wb.ForceFullCalculation = $true
